# Rename the three "Include from TRE-..." sheets to use an underscore
# after "TRE" (TRE-Rxx -> TRE_Rxx), matching the slice identifiers used
# in fr.core, and bump the Metadata "Date" value to the new publish time.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "Include from TRE-R67-TypeStru") {
        $ws.Name = "Include from TRE_R67-TypeStru"
    }
    elseif ($ws.Name -eq "Include from TRE-R04-TypeSavo") {
        $ws.Name = "Include from TRE_R04-TypeSavo"
    }
    elseif ($ws.Name -eq "Include from TRE-R288-TypePro") {
        $ws.Name = "Include from TRE_R288-TypePro"
    }
}

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-04-03T09:10:32+00:00"
